$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna4"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.434937333333333
$ws.Cells.Item(2, 8).Value = 4.304812
$ws.Cells.Item(2, 9).Value = 0.5010808920723563
$ws.Cells.Item(2, 10).Value = 0.5010808920723562
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.009849666666666666
$ws.Cells.Item(2, 14).Value = 0.029549
$ws.Cells.Item(2, 15).Value = 0.0002013876315934659
$ws.Cells.Item(2, 16).Value = 0.0002013876315934659
$ws.Cells.Item(2, 17).Value = 0.01413365442088889
$ws.Cells.Item(2, 18).Value = 0.127202889788
$ws.Cells.Item(2, 19).Value = 0.0001009114940911929
$ws.Cells.Item(2, 20).Value = 0.0001009114940911929

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna4"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.434937333333333
$ws.Cells.Item(3, 8).Value = 4.304812
$ws.Cells.Item(3, 9).Value = 0.5010808920723563
$ws.Cells.Item(3, 10).Value = 0.5010808920723562
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 46.25093466666667
$ws.Cells.Item(3, 14).Value = 138.752804
$ws.Cells.Item(3, 15).Value = 0.9456529349389956
$ws.Cells.Item(3, 16).Value = 0.9456529349389956
$ws.Cells.Item(3, 17).Value = 66.36719285476089
$ws.Cells.Item(3, 18).Value = 597.304735692848
$ws.Cells.Item(3, 19).Value = 0.4738486162300738
$ws.Cells.Item(3, 20).Value = 0.4738486162300737

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna4"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.434937333333333
$ws.Cells.Item(4, 8).Value = 4.304812
$ws.Cells.Item(4, 9).Value = 0.5010808920723563
$ws.Cells.Item(4, 10).Value = 0.5010808920723562
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.607896333333333
$ws.Cells.Item(4, 14).Value = 7.823689
$ws.Cells.Item(4, 15).Value = 0.05332140505715427
$ws.Cells.Item(4, 16).Value = 0.05332140505715428
$ws.Cells.Item(4, 17).Value = 3.742167810163111
$ws.Cells.Item(4, 18).Value = 33.679510291468
$ws.Cells.Item(4, 19).Value = 0.02671833721259031
$ws.Cells.Item(4, 20).Value = 0.02671833721259031

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna4"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.434937333333333
$ws.Cells.Item(5, 8).Value = 4.304812
$ws.Cells.Item(5, 9).Value = 0.5010808920723563
$ws.Cells.Item(5, 10).Value = 0.5010808920723562
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04031433333333333
$ws.Cells.Item(5, 14).Value = 0.120943
$ws.Cells.Item(5, 15).Value = 0.0008242723722565416
$ws.Cells.Item(5, 16).Value = 0.0008242723722565415
$ws.Cells.Item(5, 17).Value = 0.05784854196844445
$ws.Cells.Item(5, 18).Value = 0.520636877716
$ws.Cells.Item(5, 19).Value = 0.0004130271356009052
$ws.Cells.Item(5, 20).Value = 0.0004130271356009051

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna4"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9964423333333334
$ws.Cells.Item(6, 8).Value = 2.989327
$ws.Cells.Item(6, 9).Value = 0.3479582011609289
$ws.Cells.Item(6, 10).Value = 0.3479582011609288
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.009849666666666666
$ws.Cells.Item(6, 14).Value = 0.029549
$ws.Cells.Item(6, 15).Value = 0.0002013876315934659
$ws.Cells.Item(6, 16).Value = 0.0002013876315934659
$ws.Cells.Item(6, 17).Value = 0.00981462483588889
$ws.Cells.Item(6, 18).Value = 0.088331623523
$ws.Cells.Item(6, 19).Value = 0.000070074478025322254433
$ws.Cells.Item(6, 20).Value = 0.0000700744780253222273279

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna4"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9964423333333334
$ws.Cells.Item(7, 8).Value = 2.989327
$ws.Cells.Item(7, 9).Value = 0.3479582011609289
$ws.Cells.Item(7, 10).Value = 0.3479582011609288
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 46.25093466666667
$ws.Cells.Item(7, 14).Value = 138.752804
$ws.Cells.Item(7, 15).Value = 0.9456529349389956
$ws.Cells.Item(7, 16).Value = 0.9456529349389956
$ws.Cells.Item(7, 17).Value = 46.08638925810089
$ws.Cells.Item(7, 18).Value = 414.777503322908
$ws.Cells.Item(7, 19).Value = 0.3290476941639258
$ws.Cells.Item(7, 20).Value = 0.3290476941639258

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna4"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.9964423333333334
$ws.Cells.Item(8, 8).Value = 2.989327
$ws.Cells.Item(8, 9).Value = 0.3479582011609289
$ws.Cells.Item(8, 10).Value = 0.3479582011609288
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.607896333333333
$ws.Cells.Item(8, 14).Value = 7.823689
$ws.Cells.Item(8, 15).Value = 0.05332140505715427
$ws.Cells.Item(8, 16).Value = 0.05332140505715428
$ws.Cells.Item(8, 17).Value = 2.598618307478111
$ws.Cells.Item(8, 18).Value = 23.387564767303
$ws.Cells.Item(8, 19).Value = 0.01855362018706066
$ws.Cells.Item(8, 20).Value = 0.01855362018706065

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna4"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.9964423333333334
$ws.Cells.Item(9, 8).Value = 2.989327
$ws.Cells.Item(9, 9).Value = 0.3479582011609289
$ws.Cells.Item(9, 10).Value = 0.3479582011609288
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.04031433333333333
$ws.Cells.Item(9, 14).Value = 0.120943
$ws.Cells.Item(9, 15).Value = 0.0008242723722565416
$ws.Cells.Item(9, 16).Value = 0.0008242723722565415
$ws.Cells.Item(9, 17).Value = 0.04017090837344445
$ws.Cells.Item(9, 18).Value = 0.361538175361
$ws.Cells.Item(9, 19).Value = 0.0002868123319170378
$ws.Cells.Item(9, 20).Value = 0.0002868123319170376

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efna4"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.4323043333333333
$ws.Cells.Item(10, 8).Value = 1.296913
$ws.Cells.Item(10, 9).Value = 0.150960906766715
$ws.Cells.Item(10, 10).Value = 0.1509609067667149
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.009849666666666666
$ws.Cells.Item(10, 14).Value = 0.029549
$ws.Cells.Item(10, 15).Value = 0.0002013876315934659
$ws.Cells.Item(10, 16).Value = 0.0002013876315934659
$ws.Cells.Item(10, 17).Value = 0.004258053581888889
$ws.Cells.Item(10, 18).Value = 0.038322482237
$ws.Cells.Item(10, 19).Value = 0.0000304016594769507493611
$ws.Cells.Item(10, 20).Value = 0.0000304016594769507391967

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Efna4"
$ws.Cells.Item(11, 3).Value = "Epha3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.4323043333333333
$ws.Cells.Item(11, 8).Value = 1.296913
$ws.Cells.Item(11, 9).Value = 0.150960906766715
$ws.Cells.Item(11, 10).Value = 0.1509609067667149
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 46.25093466666667
$ws.Cells.Item(11, 14).Value = 138.752804
$ws.Cells.Item(11, 15).Value = 0.9456529349389956
$ws.Cells.Item(11, 16).Value = 0.9456529349389956
$ws.Cells.Item(11, 17).Value = 19.99447947711689
$ws.Cells.Item(11, 18).Value = 179.950315294052
$ws.Cells.Item(11, 19).Value = 0.1427566245449961
$ws.Cells.Item(11, 20).Value = 0.142756624544996

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Efna4"
$ws.Cells.Item(12, 3).Value = "Epha3"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.4323043333333333
$ws.Cells.Item(12, 8).Value = 1.296913
$ws.Cells.Item(12, 9).Value = 0.150960906766715
$ws.Cells.Item(12, 10).Value = 0.1509609067667149
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.607896333333333
$ws.Cells.Item(12, 14).Value = 7.823689
$ws.Cells.Item(12, 15).Value = 0.05332140505715427
$ws.Cells.Item(12, 16).Value = 0.05332140505715428
$ws.Cells.Item(12, 17).Value = 1.127404885784111
$ws.Cells.Item(12, 18).Value = 10.146643972057
$ws.Cells.Item(12, 19).Value = 0.00804944765750331
$ws.Cells.Item(12, 20).Value = 0.008049447657503308

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Efna4"
$ws.Cells.Item(13, 3).Value = "Epha3"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.4323043333333333
$ws.Cells.Item(13, 8).Value = 1.296913
$ws.Cells.Item(13, 9).Value = 0.150960906766715
$ws.Cells.Item(13, 10).Value = 0.1509609067667149
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.04031433333333333
$ws.Cells.Item(13, 14).Value = 0.120943
$ws.Cells.Item(13, 15).Value = 0.0008242723722565416
$ws.Cells.Item(13, 16).Value = 0.0008242723722565415
$ws.Cells.Item(13, 17).Value = 0.01742806099544444
$ws.Cells.Item(13, 18).Value = 0.156852548959
$ws.Cells.Item(13, 19).Value = 0.0001244329047385988
$ws.Cells.Item(13, 20).Value = 0.0001244329047385987
